# Apply the cryptos.xlsx price/volume refresh described in the commit.
# Source rows are plain text cells (inlineStr): coin name (B), link (C),
# price (D) and 1h volume % (E). A handful of D-column values look like
# plain decimals ("0.700", "73.80", ...) which Excel would otherwise
# auto-convert to numbers (dropping trailing zeros); those are written
# with a leading apostrophe to force text, then restyled back to Normal
# so no stray cell-level formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.129.64"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "1.903.19"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").Value = "'253.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.23%  "
$ws.Range("D6").Value = "'0.700"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.85%  "
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("D8").Value = "'41.39"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.69%  "
$ws.Range("D9").Value = "'0.357"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.37%  "
$ws.Range("D10").Value = "'52.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("E11").Value = "  +5.42%  "
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").Value = "'13.14"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.69%  "
$ws.Range("D14").Value = "2.180.92"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("E15").Value = "  +4.57%  "
$ws.Range("D16").Value = "'5.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.40%  "
$ws.Range("D17").Value = "1.899.50"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "35.136.86"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").Value = "'73.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.79%  "
$ws.Range("D20").Value = "0.0₃0841"
$ws.Range("E20").Value = "  +3.32%  "
$ws.Range("D21").Value = "'242.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("D22").Value = "'13.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.11%  "
$ws.Range("E23").Value = "  +6.03%  "
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("E25").Value = "  +6.41%  "
$ws.Range("D26").Value = "'2.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'167.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("D29").Value = "'18.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.78%  "
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").Value = "4.128.19"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("E32").Value = "  +7.03%  "
$ws.Range("D33").Value = "'2.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.06%  "
$ws.Range("D34").Value = "'4.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.76%  "
$ws.Range("E35").Value = "  +9.03%  "
$ws.Range("E36").Value = "  +4.22%  "
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("E38").Value = "  -6.98%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").Value = "'99.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.52%  "
$ws.Range("D41").Value = "'17.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.20%  "
$ws.Range("E42").Value = "  +4.19%  "
$ws.Range("D43").Value = "'1.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.76%  "
$ws.Range("D44").Value = "'0.0654"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.88%  "
$ws.Range("D45").Value = "'2.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").Value = "1.306.13"
$ws.Range("E46").Value = "  -3.16%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D48").Value = "'2.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("B49").Value = "Gas"
$ws.Range("C49").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D49").Value = "'12.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.59%  "
$ws.Range("E50").Value = "  +2.45%  "
$ws.Range("E51").Value = "  +7.19%  "
